$wb = $excel.ActiveWorkbook

# --- Image_Alt_Attribute_Test sheet updates ---
$ws3 = $wb.Worksheets.Item("Image_Alt_Attribute_Test")
$ws3.Range("E63").Value = "Alt attribute present: Gran oferta para esta noche"
$ws3.Range("E65").Value = "Alt attribute present: Reservado 7 veces en las últimas 24 horas"
$ws3.Range("E67").Value = "Alt attribute present: Casa mejor calificada en Porto"
$ws3.Range("E69").Value = "Alt attribute present: Casa excepcional cerca de Porto"
$ws3.Range("E71").Value = "Alt attribute present: Alquiler de viviendas de lujo a Hotala™ 99"
$ws3.Range("E73").Value = "Alt attribute present: Casa ecológica en Porto"

# --- URL_Status_Code_Test sheet updates (rows 2-136 reshuffled) ---
$ws4 = $wb.Worksheets.Item("URL_Status_Code_Test")
$arr4 = New-Object 'object[,]' 135,5
$arr4[0,0] = "https://www.alojamiento.io/all/andorra"
$arr4[0,1] = "URL Status Code Test"
$arr4[0,2] = 403
$arr4[0,3] = "passed"
$arr4[0,4] = "Status code 403 (Forbidden)"
$arr4[1,0] = "https://www.alojamiento.io/all/greece"
$arr4[1,1] = "URL Status Code Test"
$arr4[1,2] = 403
$arr4[1,3] = "passed"
$arr4[1,4] = "Status code 403 (Forbidden)"
$arr4[2,0] = "https://www.alojamiento.io/near-me?all=true"
$arr4[2,1] = "URL Status Code Test"
$arr4[2,2] = 403
$arr4[2,3] = "passed"
$arr4[2,4] = "Status code 403 (Forbidden)"
$arr4[3,0] = "https://www.alojamiento.io/all/portugal/porto-district/porto/paranhos/hotels"
$arr4[3,1] = "URL Status Code Test"
$arr4[3,2] = 403
$arr4[3,3] = "passed"
$arr4[3,4] = "Status code 403 (Forbidden)"
$arr4[4,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-2134335&guests=2&search_string=Porto,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=2-BC-2134335&published=true&dest_id=2134335&hero=BC-1935047&owner_id=2134335&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2023-09-04T08%3A44%3A02.076302%2B00%3A00"
$arr4[4,1] = "URL Status Code Test"
$arr4[4,2] = 403
$arr4[4,3] = "passed"
$arr4[4,4] = "Status code 403 (Forbidden)"
$arr4[5,0] = "https://www.alojamiento.io/all/argentina/mendoza"
$arr4[5,1] = "URL Status Code Test"
$arr4[5,2] = 403
$arr4[5,3] = "passed"
$arr4[5,4] = "Status code 403 (Forbidden)"
$arr4[6,0] = "https://www.alojamiento.io/property/romantic-luxury-house-w-captivating-outdoor-patio/BC-10659330"
$arr4[6,1] = "URL Status Code Test"
$arr4[6,2] = 403
$arr4[6,3] = "passed"
$arr4[6,4] = "Status code 403 (Forbidden)"
$arr4[7,0] = "https://www.alojamiento.io/all/honduras/bay-islands/roatan"
$arr4[7,1] = "URL Status Code Test"
$arr4[7,2] = 403
$arr4[7,3] = "passed"
$arr4[7,4] = "Status code 403 (Forbidden)"
$arr4[8,0] = "https://www.alojamiento.io/about-us"
$arr4[8,1] = "URL Status Code Test"
$arr4[8,2] = 403
$arr4[8,3] = "passed"
$arr4[8,4] = "Status code 403 (Forbidden)"
$arr4[9,0] = "https://www.alojamiento.io/all/spain/catalonia/costa-brava"
$arr4[9,1] = "URL Status Code Test"
$arr4[9,2] = 403
$arr4[9,3] = "passed"
$arr4[9,4] = "Status code 403 (Forbidden)"
$arr4[10,0] = "https://www.alojamiento.io/redirect-partner?feed=12&property_id=HA-6166825172&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=9-HA-6166825172&published=true&dest_id=18971561&hero=BC-1935047&owner_id=18971561&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&epc=c999&upat=2024-12-04T10%3A03%3A04.951386%2B00%3A00"
$arr4[10,1] = "URL Status Code Test"
$arr4[10,2] = 403
$arr4[10,3] = "passed"
$arr4[10,4] = "Status code 403 (Forbidden)"
$arr4[11,0] = "https://www.alojamiento.io/all/mexico/quintana-roo/cancun"
$arr4[11,1] = "URL Status Code Test"
$arr4[11,2] = 403
$arr4[11,3] = "passed"
$arr4[11,4] = "Status code 403 (Forbidden)"
$arr4[12,0] = "https://www.alojamiento.io/all/spain/community-of-madrid/madrid"
$arr4[12,1] = "URL Status Code Test"
$arr4[12,2] = 403
$arr4[12,3] = "passed"
$arr4[12,4] = "Status code 403 (Forbidden)"
$arr4[13,0] = "https://www.alojamiento.io/all/spain/catalonia"
$arr4[13,1] = "URL Status Code Test"
$arr4[13,2] = 403
$arr4[13,3] = "passed"
$arr4[13,4] = "Status code 403 (Forbidden)"
$arr4[14,0] = "https://www.alojamiento.io/all/italy"
$arr4[14,1] = "URL Status Code Test"
$arr4[14,2] = 403
$arr4[14,3] = "passed"
$arr4[14,4] = "Status code 403 (Forbidden)"
$arr4[15,0] = "https://www.alojamiento.io/site-map"
$arr4[15,1] = "URL Status Code Test"
$arr4[15,2] = 403
$arr4[15,3] = "passed"
$arr4[15,4] = "Status code 403 (Forbidden)"
$arr4[16,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=placeholder5&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=4-placeholder5&hero=BC-1935047&order=upsort_bh&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340"
$arr4[16,1] = "URL Status Code Test"
$arr4[16,2] = 403
$arr4[16,3] = "passed"
$arr4[16,4] = "Status code 403 (Forbidden)"
$arr4[17,0] = "https://www.alojamiento.io/property/cozy-apartment-invicta-city-opo/BC-11534173"
$arr4[17,1] = "URL Status Code Test"
$arr4[17,2] = 403
$arr4[17,3] = "passed"
$arr4[17,4] = "Status code 403 (Forbidden)"
$arr4[18,0] = "https://www.alojamiento.io/all/panama/panama/panama-city"
$arr4[18,1] = "URL Status Code Test"
$arr4[18,2] = 403
$arr4[18,3] = "passed"
$arr4[18,4] = "Status code 403 (Forbidden)"
$arr4[19,0] = "https://www.alojamiento.io/privacy-policy#site-cookie-policy"
$arr4[19,1] = "URL Status Code Test"
$arr4[19,2] = 403
$arr4[19,3] = "passed"
$arr4[19,4] = "Status code 403 (Forbidden)"
$arr4[20,0] = "https://www.alojamiento.io/all/guatemala/sacatepequez/antigua-guatemala"
$arr4[20,1] = "URL Status Code Test"
$arr4[20,2] = 403
$arr4[20,3] = "passed"
$arr4[20,4] = "Status code 403 (Forbidden)"
$arr4[21,0] = "https://www.alojamiento.io/place-to-stay"
$arr4[21,1] = "URL Status Code Test"
$arr4[21,2] = 403
$arr4[21,3] = "passed"
$arr4[21,4] = "Status code 403 (Forbidden)"
$arr4[22,0] = "https://www.alojamiento.io/property/urban-views-bright-apt-inspired-by-pal%c3%a1cio-cristal/HA-61611682440"
$arr4[22,1] = "URL Status Code Test"
$arr4[22,2] = 403
$arr4[22,3] = "passed"
$arr4[22,4] = "Status code 403 (Forbidden)"
$arr4[23,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-11172109&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=13-BC-11172109&published=true&dest_id=11172109&hero=BC-1935047&owner_id=11172109&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-12-05T15%3A02%3A53.704870%2B00%3A00"
$arr4[23,1] = "URL Status Code Test"
$arr4[23,2] = 403
$arr4[23,3] = "passed"
$arr4[23,4] = "Status code 403 (Forbidden)"
$arr4[24,0] = "https://www.alojamiento.io/all/cuba/province-of-havana"
$arr4[24,1] = "URL Status Code Test"
$arr4[24,2] = 403
$arr4[24,3] = "passed"
$arr4[24,4] = "Status code 403 (Forbidden)"
$arr4[25,0] = "https://www.alojamiento.io/addalisting"
$arr4[25,1] = "URL Status Code Test"
$arr4[25,2] = 403
$arr4[25,3] = "passed"
$arr4[25,4] = "Status code 403 (Forbidden)"
$arr4[26,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-11534173&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=10-BC-11534173&published=true&dest_id=11534173&hero=BC-1935047&owner_id=11534173&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-05-12T13%3A17%3A41.359928%2B00%3A00"
$arr4[26,1] = "URL Status Code Test"
$arr4[26,2] = 403
$arr4[26,3] = "passed"
$arr4[26,4] = "Status code 403 (Forbidden)"
$arr4[27,0] = "https://www.alojamiento.io/all/argentina/buenos-aires"
$arr4[27,1] = "URL Status Code Test"
$arr4[27,2] = 403
$arr4[27,3] = "passed"
$arr4[27,4] = "Status code 403 (Forbidden)"
$arr4[28,0] = "https://www.alojamiento.io/property/campo-lindo-apartment/BC-1935047"
$arr4[28,1] = "URL Status Code Test"
$arr4[28,2] = 403
$arr4[28,3] = "passed"
$arr4[28,4] = "Status code 403 (Forbidden)"
$arr4[29,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-10998286&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=11-BC-10998286&published=true&dest_id=10998286&hero=BC-1935047&owner_id=10998286&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-12-05T07%3A15%3A38.030190%2B00%3A00"
$arr4[29,1] = "URL Status Code Test"
$arr4[29,2] = 403
$arr4[29,3] = "passed"
$arr4[29,4] = "Status code 403 (Forbidden)"
$arr4[30,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-10599206&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=14-BC-10599206&published=true&dest_id=10599206&hero=BC-1935047&owner_id=10599206&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-07-16T03%3A31%3A23.472631%2B00%3A00"
$arr4[30,1] = "URL Status Code Test"
$arr4[30,2] = 403
$arr4[30,3] = "passed"
$arr4[30,4] = "Status code 403 (Forbidden)"
$arr4[31,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=placeholder3&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=2-placeholder3&hero=BC-1935047&order=upsort_bh&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340"
$arr4[31,1] = "URL Status Code Test"
$arr4[31,2] = 403
$arr4[31,3] = "passed"
$arr4[31,4] = "Status code 403 (Forbidden)"
$arr4[32,0] = "https://x.com/StaysTravel"
$arr4[32,1] = "URL Status Code Test"
$arr4[32,2] = 200
$arr4[32,3] = "passed"
$arr4[32,4] = "Status code 200 (OK)"
$arr4[33,0] = "https://www.alojamiento.io/all/mexico"
$arr4[33,1] = "URL Status Code Test"
$arr4[33,2] = 403
$arr4[33,3] = "passed"
$arr4[33,4] = "Status code 403 (Forbidden)"
$arr4[34,0] = "https://www.alojamiento.io/refine?search=Paranhos%2c%20Porto%2c%20Porto%20District%2c%20Portugal"
$arr4[34,1] = "URL Status Code Test"
$arr4[34,2] = 403
$arr4[34,3] = "passed"
$arr4[34,4] = "Status code 403 (Forbidden)"
$arr4[35,0] = "https://www.alojamiento.io/all/mexico/guanajuato/san-miguel-de-allende"
$arr4[35,1] = "URL Status Code Test"
$arr4[35,2] = 403
$arr4[35,3] = "passed"
$arr4[35,4] = "Status code 403 (Forbidden)"
$arr4[36,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-11099494&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=17-BC-11099494&published=true&dest_id=11099494&hero=BC-1935047&owner_id=11099494&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-12-01T07%3A14%3A52.613135%2B00%3A00"
$arr4[36,1] = "URL Status Code Test"
$arr4[36,2] = 403
$arr4[36,3] = "passed"
$arr4[36,4] = "Status code 403 (Forbidden)"
$arr4[37,0] = "https://www.alojamiento.io/all/brazil/south-region/florianopolis"
$arr4[37,1] = "URL Status Code Test"
$arr4[37,2] = 403
$arr4[37,3] = "passed"
$arr4[37,4] = "Status code 403 (Forbidden)"
$arr4[38,0] = "https://www.alojamiento.io/property/vila-cam%c3%a9lia/BC-11099494"
$arr4[38,1] = "URL Status Code Test"
$arr4[38,2] = 403
$arr4[38,3] = "passed"
$arr4[38,4] = "Status code 403 (Forbidden)"
$arr4[39,0] = "https://www.alojamiento.io/redirect-partner?feed=12&property_id=HA-6118850658&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=5-HA-6118850658&published=true&dest_id=29545296&hero=BC-1935047&owner_id=29545296&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&epc=c999&upat=2024-12-04T11%3A48%3A37.454906%2B00%3A00"
$arr4[39,1] = "URL Status Code Test"
$arr4[39,2] = 403
$arr4[39,3] = "passed"
$arr4[39,4] = "Status code 403 (Forbidden)"
$arr4[40,0] = "https://www.alojamiento.io/all/costa-rica/san-jose"
$arr4[40,1] = "URL Status Code Test"
$arr4[40,2] = 403
$arr4[40,3] = "passed"
$arr4[40,4] = "Status code 403 (Forbidden)"
$arr4[41,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-11242224&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=8-BC-11242224&published=true&dest_id=11242224&hero=BC-1935047&owner_id=11242224&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-12-05T17%3A21%3A13.957523%2B00%3A00"
$arr4[41,1] = "URL Status Code Test"
$arr4[41,2] = 403
$arr4[41,3] = "passed"
$arr4[41,4] = "Status code 403 (Forbidden)"
$arr4[42,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-4649516&guests=2&search_string=Porto,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=19-BC-4649516&published=true&dest_id=4649516&hero=BC-1935047&owner_id=4649516&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-08-24T10%3A49%3A52.209651%2B00%3A00"
$arr4[42,1] = "URL Status Code Test"
$arr4[42,2] = 403
$arr4[42,3] = "passed"
$arr4[42,4] = "Status code 403 (Forbidden)"
$arr4[43,0] = "https://www.alojamiento.io/property/modern-apartment-with-balcony-and-free-parking/BC-11172109"
$arr4[43,1] = "URL Status Code Test"
$arr4[43,2] = 403
$arr4[43,3] = "passed"
$arr4[43,4] = "Status code 403 (Forbidden)"
$arr4[44,0] = "https://www.alojamiento.io/all/morocco"
$arr4[44,1] = "URL Status Code Test"
$arr4[44,2] = 403
$arr4[44,3] = "passed"
$arr4[44,4] = "Status code 403 (Forbidden)"
$arr4[45,0] = "https://www.alojamiento.io/all/belize/belize-district/belize-city"
$arr4[45,1] = "URL Status Code Test"
$arr4[45,2] = 403
$arr4[45,3] = "passed"
$arr4[45,4] = "Status code 403 (Forbidden)"
$arr4[46,0] = "https://www.alojamiento.io/all/portugal/porto-district/porto/paranhos/vacation-rentals"
$arr4[46,1] = "URL Status Code Test"
$arr4[46,2] = 403
$arr4[46,3] = "passed"
$arr4[46,4] = "Status code 403 (Forbidden)"
$arr4[47,0] = "https://www.travelai.com/"
$arr4[47,1] = "URL Status Code Test"
$arr4[47,2] = 200
$arr4[47,3] = "passed"
$arr4[47,4] = "Status code 200 (OK)"
$arr4[48,0] = "https://www.alojamiento.io/all/portugal/porto-district/porto/paranhos/resorts"
$arr4[48,1] = "URL Status Code Test"
$arr4[48,2] = 403
$arr4[48,3] = "passed"
$arr4[48,4] = "Status code 403 (Forbidden)"
$arr4[49,0] = "https://www.alojamiento.io/all/colombia/medellin"
$arr4[49,1] = "URL Status Code Test"
$arr4[49,2] = 403
$arr4[49,3] = "passed"
$arr4[49,4] = "Status code 403 (Forbidden)"
$arr4[50,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-12847410&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=23-BC-12847410&published=true&dest_id=12847410&hero=BC-1935047&owner_id=12847410&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-12-02T08%3A18%3A47.203235%2B00%3A00"
$arr4[50,1] = "URL Status Code Test"
$arr4[50,2] = 403
$arr4[50,3] = "passed"
$arr4[50,4] = "Status code 403 (Forbidden)"
$arr4[51,0] = "https://www.alojamiento.io/property/poetikblue-by-we-do-living/BC-11346246"
$arr4[51,1] = "URL Status Code Test"
$arr4[51,2] = 403
$arr4[51,3] = "passed"
$arr4[51,4] = "Status code 403 (Forbidden)"
$arr4[52,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=placeholder2&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=1-placeholder2&hero=BC-1935047&order=upsort_bh&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340"
$arr4[52,1] = "URL Status Code Test"
$arr4[52,2] = 403
$arr4[52,3] = "passed"
$arr4[52,4] = "Status code 403 (Forbidden)"
$arr4[53,0] = "https://www.alojamiento.io/all/mexico/baja-california-sur/cabo-san-lucas"
$arr4[53,1] = "URL Status Code Test"
$arr4[53,2] = 403
$arr4[53,3] = "passed"
$arr4[53,4] = "Status code 403 (Forbidden)"
$arr4[54,0] = "https://www.alojamiento.io/all/spain/canary-islands"
$arr4[54,1] = "URL Status Code Test"
$arr4[54,2] = 403
$arr4[54,3] = "passed"
$arr4[54,4] = "Status code 403 (Forbidden)"
$arr4[55,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-12847446&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=15-BC-12847446&published=true&dest_id=12847446&hero=BC-1935047&owner_id=12847446&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-12-02T08%3A19%3A35.025518%2B00%3A00"
$arr4[55,1] = "URL Status Code Test"
$arr4[55,2] = 403
$arr4[55,3] = "passed"
$arr4[55,4] = "Status code 403 (Forbidden)"
$arr4[56,0] = "https://www.alojamiento.io/property/s-jo%c3%a3o-porto-apartment/BC-2202420"
$arr4[56,1] = "URL Status Code Test"
$arr4[56,2] = 403
$arr4[56,3] = "passed"
$arr4[56,4] = "Status code 403 (Forbidden)"
$arr4[57,0] = "https://www.alojamiento.io/all/spain/navarre"
$arr4[57,1] = "URL Status Code Test"
$arr4[57,2] = 403
$arr4[57,3] = "passed"
$arr4[57,4] = "Status code 403 (Forbidden)"
$arr4[58,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=placeholder6&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=5-placeholder6&hero=BC-1935047&order=upsort_bh&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340"
$arr4[58,1] = "URL Status Code Test"
$arr4[58,2] = 403
$arr4[58,3] = "passed"
$arr4[58,4] = "Status code 403 (Forbidden)"
$arr4[59,0] = "https://www.alojamiento.io/all/spain/castile-and-leon"
$arr4[59,1] = "URL Status Code Test"
$arr4[59,2] = 403
$arr4[59,3] = "passed"
$arr4[59,4] = "Status code 403 (Forbidden)"
$arr4[60,0] = "https://www.alojamiento.io/all/mexico/jalisco/puerto-vallarta"
$arr4[60,1] = "URL Status Code Test"
$arr4[60,2] = 403
$arr4[60,3] = "passed"
$arr4[60,4] = "Status code 403 (Forbidden)"
$arr4[61,0] = "https://www.alojamiento.io/all/portugal/porto-district/porto/paranhos/villas"
$arr4[61,1] = "URL Status Code Test"
$arr4[61,2] = 403
$arr4[61,3] = "passed"
$arr4[61,4] = "Status code 403 (Forbidden)"
$arr4[62,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-11105844&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=22-BC-11105844&published=true&dest_id=11105844&hero=BC-1935047&owner_id=11105844&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-12-05T11%3A25%3A32.784787%2B00%3A00"
$arr4[62,1] = "URL Status Code Test"
$arr4[62,2] = 403
$arr4[62,3] = "passed"
$arr4[62,4] = "Status code 403 (Forbidden)"
$arr4[63,0] = "https://www.alojamiento.io/all/france"
$arr4[63,1] = "URL Status Code Test"
$arr4[63,2] = 403
$arr4[63,3] = "passed"
$arr4[63,4] = "Status code 403 (Forbidden)"
$arr4[64,0] = "https://www.alojamiento.io/all/colombia/magdalena/santa-marta"
$arr4[64,1] = "URL Status Code Test"
$arr4[64,2] = 403
$arr4[64,3] = "passed"
$arr4[64,4] = "Status code 403 (Forbidden)"
$arr4[65,0] = "https://www.alojamiento.io/property/ac-house/BC-2807903"
$arr4[65,1] = "URL Status Code Test"
$arr4[65,2] = 403
$arr4[65,3] = "passed"
$arr4[65,4] = "Status code 403 (Forbidden)"
$arr4[66,0] = "https://www.alojamiento.io/all/usa/wyoming/basin"
$arr4[66,1] = "URL Status Code Test"
$arr4[66,2] = 403
$arr4[66,3] = "passed"
$arr4[66,4] = "Status code 403 (Forbidden)"
$arr4[67,0] = "https://www.alojamiento.io/all/brazil/southeast-region/rio-de-janeiro"
$arr4[67,1] = "URL Status Code Test"
$arr4[67,2] = 403
$arr4[67,3] = "passed"
$arr4[67,4] = "Status code 403 (Forbidden)"
$arr4[68,0] = "https://www.alojamiento.io/property/lv-premier-marques-mq2-balc%c3%b3n-aire-acondicionado-vistas/HA-6118850658"
$arr4[68,1] = "URL Status Code Test"
$arr4[68,2] = 403
$arr4[68,3] = "passed"
$arr4[68,4] = "Status code 403 (Forbidden)"
$arr4[69,0] = "https://www.alojamiento.io/redirect-partner?feed=12&property_id=HA-61611682438&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=3-HA-61611682438&published=true&dest_id=109825542&hero=BC-1935047&owner_id=109825542&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&epc=c001&upat=2024-12-05T01%3A44%3A21.373667%2B00%3A00"
$arr4[69,1] = "URL Status Code Test"
$arr4[69,2] = 403
$arr4[69,3] = "passed"
$arr4[69,4] = "Status code 403 (Forbidden)"
$arr4[70,0] = "https://www.alojamiento.io/all/peru/cusco/cusco"
$arr4[70,1] = "URL Status Code Test"
$arr4[70,2] = 403
$arr4[70,3] = "passed"
$arr4[70,4] = "Status code 403 (Forbidden)"
$arr4[71,0] = "https://www.alojamiento.io/property/portogali-guesthouse/BC-5813166"
$arr4[71,1] = "URL Status Code Test"
$arr4[71,2] = 403
$arr4[71,3] = "passed"
$arr4[71,4] = "Status code 403 (Forbidden)"
$arr4[72,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-12500411&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=20-BC-12500411&published=true&dest_id=12500411&hero=BC-1935047&owner_id=12500411&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-12-04T11%3A00%3A14.620270%2B00%3A00"
$arr4[72,1] = "URL Status Code Test"
$arr4[72,2] = 403
$arr4[72,3] = "passed"
$arr4[72,4] = "Status code 403 (Forbidden)"
$arr4[73,0] = "https://www.instagram.com/staystravel"
$arr4[73,1] = "URL Status Code Test"
$arr4[73,2] = 200
$arr4[73,3] = "passed"
$arr4[73,4] = "Status code 200 (OK)"
$arr4[74,0] = "https://www.alojamiento.io/site-terms"
$arr4[74,1] = "URL Status Code Test"
$arr4[74,2] = 403
$arr4[74,3] = "passed"
$arr4[74,4] = "Status code 403 (Forbidden)"
$arr4[75,0] = "https://www.alojamiento.io/all/chile/santiago-metropolitan/santiago"
$arr4[75,1] = "URL Status Code Test"
$arr4[75,2] = 403
$arr4[75,3] = "passed"
$arr4[75,4] = "Status code 403 (Forbidden)"
$arr4[76,0] = "https://www.alojamiento.io/all/portugal"
$arr4[76,1] = "URL Status Code Test"
$arr4[76,2] = 403
$arr4[76,3] = "passed"
$arr4[76,4] = "Status code 403 (Forbidden)"
$arr4[77,0] = "https://www.alojamiento.io/all/bolivia/la-paz"
$arr4[77,1] = "URL Status Code Test"
$arr4[77,2] = 403
$arr4[77,3] = "passed"
$arr4[77,4] = "Status code 403 (Forbidden)"
$arr4[78,0] = "https://www.alojamiento.io/all/ecuador/pichincha/quito"
$arr4[78,1] = "URL Status Code Test"
$arr4[78,2] = 403
$arr4[78,3] = "passed"
$arr4[78,4] = "Status code 403 (Forbidden)"
$arr4[79,0] = "https://www.alojamiento.io/redirect-partner?feed=12&property_id=HA-6166491756&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=7-HA-6166491756&published=true&dest_id=18976742&hero=BC-1935047&owner_id=18976742&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&epc=c999&upat=2024-12-04T10%3A04%3A28.809814%2B00%3A00"
$arr4[79,1] = "URL Status Code Test"
$arr4[79,2] = 403
$arr4[79,3] = "passed"
$arr4[79,4] = "Status code 403 (Forbidden)"
$arr4[80,0] = "https://www.alojamiento.io/all/portugal/porto-district/porto/paranhos/cottages"
$arr4[80,1] = "URL Status Code Test"
$arr4[80,2] = 403
$arr4[80,3] = "passed"
$arr4[80,4] = "Status code 403 (Forbidden)"
$arr4[81,0] = "https://www.alojamiento.io/all/spain/andalusia"
$arr4[81,1] = "URL Status Code Test"
$arr4[81,2] = 403
$arr4[81,3] = "passed"
$arr4[81,4] = "Status code 403 (Forbidden)"
$arr4[82,0] = "https://www.alojamiento.io/all/portugal/porto-district/porto/paranhos/cabins"
$arr4[82,1] = "URL Status Code Test"
$arr4[82,2] = 403
$arr4[82,3] = "passed"
$arr4[82,4] = "Status code 403 (Forbidden)"
$arr4[83,0] = "https://www.alojamiento.io/all/portugal/porto-district/porto"
$arr4[83,1] = "URL Status Code Test"
$arr4[83,2] = 403
$arr4[83,3] = "passed"
$arr4[83,4] = "Status code 403 (Forbidden)"
$arr4[84,0] = "https://www.alojamiento.io/property/rendez-vous-porto-ii-t1/HA-6166825172"
$arr4[84,1] = "URL Status Code Test"
$arr4[84,2] = 403
$arr4[84,3] = "passed"
$arr4[84,4] = "Status code 403 (Forbidden)"
$arr4[85,0] = "https://www.alojamiento.io/property/jm-alojamento-local-no-porto/BC-11242224"
$arr4[85,1] = "URL Status Code Test"
$arr4[85,2] = 403
$arr4[85,3] = "passed"
$arr4[85,4] = "Status code 403 (Forbidden)"
$arr4[86,0] = "https://www.alojamiento.io/all/argentina/rio-negro/san-carlos-de-bariloche"
$arr4[86,1] = "URL Status Code Test"
$arr4[86,2] = 403
$arr4[86,3] = "passed"
$arr4[86,4] = "Status code 403 (Forbidden)"
$arr4[87,0] = "https://www.alojamiento.io/all/spain/valencian-community/valencia-province"
$arr4[87,1] = "URL Status Code Test"
$arr4[87,2] = 403
$arr4[87,3] = "passed"
$arr4[87,4] = "Status code 403 (Forbidden)"
$arr4[88,0] = "https://www.alojamiento.io/all/ecuador/galapagos"
$arr4[88,1] = "URL Status Code Test"
$arr4[88,2] = 403
$arr4[88,3] = "passed"
$arr4[88,4] = "Status code 403 (Forbidden)"
$arr4[89,0] = "https://www.onedegreeleft.com/"
$arr4[89,1] = "URL Status Code Test"
$arr4[89,2] = 200
$arr4[89,3] = "passed"
$arr4[89,4] = "Status code 200 (OK)"
$arr4[90,0] = "https://www.alojamiento.io/redirect-partner?feed=12&property_id=HA-61611682440&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=1-HA-61611682440&published=true&dest_id=109825564&hero=BC-1935047&owner_id=109825564&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&epc=c001&upat=2024-12-05T01%3A44%3A21.373975%2B00%3A00"
$arr4[90,1] = "URL Status Code Test"
$arr4[90,2] = 403
$arr4[90,3] = "passed"
$arr4[90,4] = "Status code 403 (Forbidden)"
$arr4[91,0] = "https://www.alojamiento.io/all/peru/lima"
$arr4[91,1] = "URL Status Code Test"
$arr4[91,2] = 403
$arr4[91,3] = "passed"
$arr4[91,4] = "Status code 403 (Forbidden)"
$arr4[92,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=placeholder1&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=0-placeholder1&hero=BC-1935047&order=upsort_bh&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340"
$arr4[92,1] = "URL Status Code Test"
$arr4[92,2] = 403
$arr4[92,3] = "passed"
$arr4[92,4] = "Status code 403 (Forbidden)"
$arr4[93,0] = "https://www.facebook.com/StaysTravel"
$arr4[93,1] = "URL Status Code Test"
$arr4[93,2] = 200
$arr4[93,3] = "passed"
$arr4[93,4] = "Status code 200 (OK)"
$arr4[94,0] = "https://www.alojamiento.io/faq"
$arr4[94,1] = "URL Status Code Test"
$arr4[94,2] = 403
$arr4[94,3] = "passed"
$arr4[94,4] = "Status code 403 (Forbidden)"
$arr4[95,0] = "https://www.alojamiento.io/property/the-cale-antero-de-quental-31/BC-12847410"
$arr4[95,1] = "URL Status Code Test"
$arr4[95,2] = 403
$arr4[95,3] = "passed"
$arr4[95,4] = "Status code 403 (Forbidden)"
$arr4[96,0] = "https://www.alojamiento.io/all/switzerland"
$arr4[96,1] = "URL Status Code Test"
$arr4[96,2] = 403
$arr4[96,3] = "passed"
$arr4[96,4] = "Status code 403 (Forbidden)"
$arr4[97,0] = "https://www.alojamiento.io/property/quarto-duplo-com-casa-de-banho-exclusiva-no-porto/BC-10599206"
$arr4[97,1] = "URL Status Code Test"
$arr4[97,2] = 403
$arr4[97,3] = "passed"
$arr4[97,4] = "Status code 403 (Forbidden)"
$arr4[98,0] = "https://www.alojamiento.io/all/spain/murcia"
$arr4[98,1] = "URL Status Code Test"
$arr4[98,2] = 403
$arr4[98,3] = "passed"
$arr4[98,4] = "Status code 403 (Forbidden)"
$arr4[99,0] = "https://www.alojamiento.io/property/t1-apartamento-moderno-luz-tranquila-y-confortable-con-garaje-privado/HA-6166491756"
$arr4[99,1] = "URL Status Code Test"
$arr4[99,2] = 403
$arr4[99,3] = "passed"
$arr4[99,4] = "Status code 403 (Forbidden)"
$arr4[100,0] = "https://www.alojamiento.io/all/dominican-republic/la-altagracia/punta-cana"
$arr4[100,1] = "URL Status Code Test"
$arr4[100,2] = 403
$arr4[100,3] = "passed"
$arr4[100,4] = "Status code 403 (Forbidden)"
$arr4[101,0] = "https://www.petfriendly.io/"
$arr4[101,1] = "URL Status Code Test"
$arr4[101,2] = 403
$arr4[101,3] = "passed"
$arr4[101,4] = "Status code 403 (Forbidden)"
$arr4[102,0] = "https://www.alojamiento.io/all/portugal/porto-district/porto/paranhos"
$arr4[102,1] = "URL Status Code Test"
$arr4[102,2] = 403
$arr4[102,3] = "passed"
$arr4[102,4] = "Status code 403 (Forbidden)"
$arr4[103,0] = "https://www.alojamiento.io/all/spain/balearic-islands"
$arr4[103,1] = "URL Status Code Test"
$arr4[103,2] = 403
$arr4[103,3] = "passed"
$arr4[103,4] = "Status code 403 (Forbidden)"
$arr4[104,0] = "https://www.alojamiento.io/all/spain/cantabria"
$arr4[104,1] = "URL Status Code Test"
$arr4[104,2] = 403
$arr4[104,3] = "passed"
$arr4[104,4] = "Status code 403 (Forbidden)"
$arr4[105,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-10659330&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=12-BC-10659330&published=true&dest_id=10659330&hero=BC-1935047&owner_id=10659330&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-12-04T06%3A05%3A00.026562%2B00%3A00"
$arr4[105,1] = "URL Status Code Test"
$arr4[105,2] = 403
$arr4[105,3] = "passed"
$arr4[105,4] = "Status code 403 (Forbidden)"
$arr4[106,0] = "https://www.alojamiento.io/property/casa-s%c3%a3o-dinis/BC-2134335"
$arr4[106,1] = "URL Status Code Test"
$arr4[106,2] = 403
$arr4[106,3] = "passed"
$arr4[106,4] = "Status code 403 (Forbidden)"
$arr4[107,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-11346246&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=21-BC-11346246&published=true&dest_id=11346246&hero=BC-1935047&owner_id=11346246&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-09-22T08%3A32%3A12.548999%2B00%3A00"
$arr4[107,1] = "URL Status Code Test"
$arr4[107,2] = 403
$arr4[107,3] = "passed"
$arr4[107,4] = "Status code 403 (Forbidden)"
$arr4[108,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-1935047&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=0-BC-1935047&published=true&dest_id=1935047&hero=BC-1935047&owner_id=1935047&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2023-11-19T06%3A10%3A57.061047%2B00%3A00"
$arr4[108,1] = "URL Status Code Test"
$arr4[108,2] = 403
$arr4[108,3] = "passed"
$arr4[108,4] = "Status code 403 (Forbidden)"
$arr4[109,0] = "https://www.alojamiento.io/all/spain/galicia"
$arr4[109,1] = "URL Status Code Test"
$arr4[109,2] = 403
$arr4[109,3] = "passed"
$arr4[109,4] = "Status code 403 (Forbidden)"
$arr4[110,0] = "https://www.alojamiento.io/property/the-cale-antero-de-quental-32/BC-12847446"
$arr4[110,1] = "URL Status Code Test"
$arr4[110,2] = 403
$arr4[110,3] = "passed"
$arr4[110,4] = "Status code 403 (Forbidden)"
$arr4[111,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-2782646&guests=2&search_string=Porto,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=18-BC-2782646&published=true&dest_id=2782646&hero=BC-1935047&owner_id=2782646&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-07-06T15%3A45%3A49.219738%2B00%3A00"
$arr4[111,1] = "URL Status Code Test"
$arr4[111,2] = 403
$arr4[111,3] = "passed"
$arr4[111,4] = "Status code 403 (Forbidden)"
$arr4[112,0] = "https://www.alojamiento.io/all/colombia/bolivar/cartagena"
$arr4[112,1] = "URL Status Code Test"
$arr4[112,2] = 403
$arr4[112,3] = "passed"
$arr4[112,4] = "Status code 403 (Forbidden)"
$arr4[113,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=placeholder4&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=3-placeholder4&hero=BC-1935047&order=upsort_bh&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340"
$arr4[113,1] = "URL Status Code Test"
$arr4[113,2] = 403
$arr4[113,3] = "passed"
$arr4[113,4] = "Status code 403 (Forbidden)"
$arr4[114,0] = "https://www.alojamiento.io/all/spain/canary-islands/gran-canaria"
$arr4[114,1] = "URL Status Code Test"
$arr4[114,2] = 403
$arr4[114,3] = "passed"
$arr4[114,4] = "Status code 403 (Forbidden)"
$arr4[115,0] = "https://www.alojamiento.io/property/fred-s-house-asprela/BC-4649516"
$arr4[115,1] = "URL Status Code Test"
$arr4[115,2] = 403
$arr4[115,3] = "passed"
$arr4[115,4] = "Status code 403 (Forbidden)"
$arr4[116,0] = "https://www.alojamiento.io/all/portugal/porto-district/porto/paranhos/ski-chalets"
$arr4[116,1] = "URL Status Code Test"
$arr4[116,2] = 403
$arr4[116,3] = "passed"
$arr4[116,4] = "Status code 403 (Forbidden)"
$arr4[117,0] = "https://www.alojamiento.io/privacy-policy"
$arr4[117,1] = "URL Status Code Test"
$arr4[117,2] = 403
$arr4[117,3] = "passed"
$arr4[117,4] = "Status code 403 (Forbidden)"
$arr4[118,0] = "https://www.alojamiento.io/all/portugal/porto-district"
$arr4[118,1] = "URL Status Code Test"
$arr4[118,2] = 403
$arr4[118,3] = "passed"
$arr4[118,4] = "Status code 403 (Forbidden)"
$arr4[119,0] = "https://www.alojamiento.io/property/lovely-flat-marqu%c3%aas/BC-11105844"
$arr4[119,1] = "URL Status Code Test"
$arr4[119,2] = 403
$arr4[119,3] = "passed"
$arr4[119,4] = "Status code 403 (Forbidden)"
$arr4[120,0] = "https://www.alojamiento.io/all/uruguay/montevideo"
$arr4[120,1] = "URL Status Code Test"
$arr4[120,2] = 403
$arr4[120,3] = "passed"
$arr4[120,4] = "Status code 403 (Forbidden)"
$arr4[121,0] = "https://www.alojamiento.io/all/brazil/northeast-region/salvador"
$arr4[121,1] = "URL Status Code Test"
$arr4[121,2] = 403
$arr4[121,3] = "passed"
$arr4[121,4] = "Status code 403 (Forbidden)"
$arr4[122,0] = "https://www.alojamiento.io/all/mexico/quintana-roo/tulum"
$arr4[122,1] = "URL Status Code Test"
$arr4[122,2] = 403
$arr4[122,3] = "passed"
$arr4[122,4] = "Status code 403 (Forbidden)"
$arr4[123,0] = "https://www.alojamiento.io/all/spain/canary-islands/tenerife"
$arr4[123,1] = "URL Status Code Test"
$arr4[123,2] = 403
$arr4[123,3] = "passed"
$arr4[123,4] = "Status code 403 (Forbidden)"
$arr4[124,0] = "https://www.alojamiento.io/all/mexico/quintana-roo/playa-del-carmen"
$arr4[124,1] = "URL Status Code Test"
$arr4[124,2] = 403
$arr4[124,3] = "passed"
$arr4[124,4] = "Status code 403 (Forbidden)"
$arr4[125,0] = "https://www.alojamiento.io/property/urban-views-petfrien-ap-inspired-by-torre-cl%c3%a9rigos/HA-61611682438"
$arr4[125,1] = "URL Status Code Test"
$arr4[125,2] = 403
$arr4[125,3] = "passed"
$arr4[125,4] = "Status code 403 (Forbidden)"
$arr4[126,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-5813166&guests=2&search_string=Paranhos,%20Porto,%20Porto%20District,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=6-BC-5813166&published=true&dest_id=5813166&hero=BC-1935047&owner_id=5813166&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2023-08-16T12%3A45%3A47.750380%2B00%3A00"
$arr4[126,1] = "URL Status Code Test"
$arr4[126,2] = 403
$arr4[126,3] = "passed"
$arr4[126,4] = "Status code 403 (Forbidden)"
$arr4[127,0] = "https://www.alojamiento.io/property/casa-velha-guesthouse/BC-2782646"
$arr4[127,1] = "URL Status Code Test"
$arr4[127,2] = 403
$arr4[127,3] = "passed"
$arr4[127,4] = "Status code 403 (Forbidden)"
$arr4[128,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-2202420&guests=2&search_string=Porto,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=4-BC-2202420&published=true&dest_id=2202420&hero=BC-1935047&owner_id=2202420&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2023-12-06T11%3A00%3A04.420388%2B00%3A00"
$arr4[128,1] = "URL Status Code Test"
$arr4[128,2] = 403
$arr4[128,3] = "passed"
$arr4[128,4] = "Status code 403 (Forbidden)"
$arr4[129,0] = "https://www.alojamiento.io/all/puerto-rico/san-juan"
$arr4[129,1] = "URL Status Code Test"
$arr4[129,2] = 403
$arr4[129,3] = "passed"
$arr4[129,4] = "Status code 403 (Forbidden)"
$arr4[130,0] = "https://www.alojamiento.io/redirect-partner?feed=11&property_id=BC-2807903&guests=2&search_string=Porto,%20Portugal&referrer_page=hybrid&menu_id=1733847445262&referral_id=16-BC-2807903&published=true&dest_id=2807903&hero=BC-1935047&owner_id=2807903&sqs=hybrid-default&property_country=PT&at=End-of-Result%20Ad&eplId=6336340&upat=2024-11-19T09%3A57%3A24.866726%2B00%3A00"
$arr4[130,1] = "URL Status Code Test"
$arr4[130,2] = 403
$arr4[130,3] = "passed"
$arr4[130,4] = "Status code 403 (Forbidden)"
$arr4[131,0] = "https://www.alojamiento.io/all/spain/andalusia/malaga"
$arr4[131,1] = "URL Status Code Test"
$arr4[131,2] = 403
$arr4[131,3] = "passed"
$arr4[131,4] = "Status code 403 (Forbidden)"
$arr4[132,0] = "https://www.alojamiento.io/property/lux-in-porto/BC-10998286"
$arr4[132,1] = "URL Status Code Test"
$arr4[132,2] = 403
$arr4[132,3] = "passed"
$arr4[132,4] = "Status code 403 (Forbidden)"
$arr4[133,0] = "https://www.alojamiento.io/property/renovado-com-ac-2-wc-completos-terra%c3%a7o-e-garagem/BC-12500411"
$arr4[133,1] = "URL Status Code Test"
$arr4[133,2] = 403
$arr4[133,3] = "passed"
$arr4[133,4] = "Status code 403 (Forbidden)"
$arr4[134,0] = "https://www.alojamiento.io/"
$arr4[134,1] = "URL Status Code Test"
$arr4[134,2] = 403
$arr4[134,3] = "passed"
$arr4[134,4] = "Status code 403 (Forbidden)"
$ws4.Range("A2:E136").Value = $arr4

Write-Host "Done"
